# daily auto push: 2026-01-22 18:51 UTC
# Insert a new data row just above the 2026/12/29 block (new row 695) so
# that a "2026/01/22" / 木 / 23 / 201 reading is recorded, pushing every
# subsequent row down by one (old row 736 -> new row 737).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 695 (and everything below it) down by one row.
$ws.Rows("695:695").Insert()

# Write the new row's values. Force column A to Text format first so the
# date-like string "2026/01/22" is stored as a literal string (matching
# the rest of the sheet) instead of being auto-converted to a date serial.
$ws.Range("A695").NumberFormat = "@"
$ws.Range("A695").Value = "2026/01/22"
$ws.Range("B695").Value = "木"
$ws.Range("C695").Value = 23
$ws.Range("D695").Value = 201

# Re-apply the plain (unstyled) look used by every other data row, since
# the temporary Text number format above would otherwise leave a stray
# style on A695.
$ws.Range("A695").Style = $ws.Range("A2").Style
